$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.05546469241380692, 0.9857710003852844, 0.02120971865952015, 0.9969012141227722),
    @(0.008792893961071968, 0.9985515475273132, 0.01113342307507992, 0.9974480271339417),
    @(0.004834053106606007, 0.9988753199577332, 0.009319002740085125, 0.9983594417572021),
    @(0.002686832100152969, 0.9991990923881531, 0.0302970577031374, 0.9925264120101929),
    @(0.001079651061445475, 0.9996762275695801, 0.004310355987399817, 0.9990885853767395),
    @(0.001054073916748166, 0.9997784495353699, 0.006501364521682262, 0.9992708563804626),
    @(0.001378633547574282, 0.9996932744979858, 0.002667276421561837, 0.9998177289962769),
    @(0.0008272122358903289, 0.999914824962616, 0.003131251782178879, 0.9996354579925537),
    @(0.0006469819927588105, 0.9998466372489929, 0.004196631256490946, 0.9996354579925537),
    @(0.0006968732341192663, 0.9997614026069641, 0.003991793841123581, 0.9994531273841858),
    @(0.0005836400086991489, 0.9998636841773987, 0.003008679719641805, 0.9998177289962769),
    @(0.0005657775327563286, 0.9998295903205872, 0.002799747278913856, 0.9998177289962769),
    @(0.0003733697230927646, 0.999914824962616, 0.004434531088918447, 0.9992708563804626),
    @(0.0005505099543370306, 0.9998295903205872, 0.00376470317132771, 0.9998177289962769),
    @(0.0001653241779422387, 0.9999659061431885, 0.003182498039677739, 0.9998177289962769),
    @(0.0001914775348268449, 0.999931812286377, 0.002444100100547075, 0.9994531273841858),
    @(0.0006572164129465818, 0.9998295903205872, 0.003061622846871614, 0.9998177289962769),
    @(0.0001622388954274356, 0.999931812286377, 0.004595072474330664, 0.9996354579925537),
    @(0.00004482537769945338, 1, 0.004683293402194977, 0.9998177289962769),
    @(0.000906155735719949, 0.9998125433921814, 0.006233714520931244, 0.9994531273841858),
    @(0.0001905623648781329, 0.999931812286377, 0.005516344681382179, 0.9996354579925537),
    @(0.0001762445026542991, 0.999931812286377, 0.005234994925558567, 0.9998177289962769),
    @(0.0003100775938946754, 0.999914824962616, 0.005051994230598211, 0.9994531273841858),
    @(0.00007663914584554732, 0.9999829530715942, 0.00594202009961009, 0.9992708563804626),
    @(0.001122172921895981, 0.9998977780342102, 0.007453171536326408, 0.9992708563804626),
    @(0.0002121882571373135, 0.9999488592147827, 0.005154968705028296, 0.9998177289962769),
    @(0.00002064303589577321, 1, 0.005261867307126522, 0.9998177289962769),
    @(0.00005666606375598349, 0.9999829530715942, 0.004669816233217716, 0.9998177289962769),
    @(0.000384041981305927, 0.999914824962616, 0.00964447483420372, 0.9989063143730164),
    @(0.0009051908855326474, 0.9998977780342102, 0.004058398306369781, 0.9998177289962769),
    @(0.00009731172031024471, 0.9999659061431885, 0.004479608964174986, 0.9998177289962769),
    @(0.0004871474229730666, 0.9998807311058044, 0.003534258343279362, 0.9998177289962769),
    @(0.0001716290134936571, 0.9999488592147827, 0.004961859900504351, 0.9998177289962769),
    @(0.000564004760235548, 0.9998977780342102, 0.005187858827412128, 0.9998177289962769),
    @(0.0001657463290030137, 0.999914824962616, 0.004568100441247225, 0.9998177289962769),
    @(0.000376501091523096, 0.9999488592147827, 0.003094392595812678, 0.9998177289962769),
    @(0.000007370269941020524, 1, 0.003056098939850926, 0.9998177289962769),
    @(0.00001944142604770605, 1, 0.003786645131185651, 0.9998177289962769),
    @(0.0001192198251374066, 0.9999829530715942, 0.002626066096127033, 0.9998177289962769),
    @(0.0003875205875374377, 0.999914824962616, 0.004368608351796865, 0.9998177289962769),
    @(0.0004514518950600177, 0.9998977780342102, 0.003108926117420197, 0.9998177289962769),
    @(0.00008368898852495477, 0.9999659061431885, 0.00434858025982976, 0.9998177289962769),
    @(0.00002437766124785412, 1, 0.004876092541962862, 0.9998177289962769),
    @(0.0001684932212810963, 0.9999659061431885, 0.006999948062002659, 0.9992708563804626),
    @(0.0005449990276247263, 0.999931812286377, 0.005144125781953335, 0.9998177289962769),
    @(0.0002141451550414786, 0.9999829530715942, 0.006610941141843796, 0.9994531273841858),
    @(0.0003499342419672757, 0.999931812286377, 0.007483762223273516, 0.9994531273841858),
    @(0.001280618598684669, 0.9998636841773987, 0.005615551024675369, 0.9998177289962769),
    @(0.0001882129727164283, 0.9999659061431885, 0.005170683842152357, 0.9998177289962769),
    @(0.00002549076998548117, 0.9999829530715942, 0.00534038245677948, 0.9998177289962769)

)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
